$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.650.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.06%  "

$ws.Range("D3").Value = "'1.801.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.62%  "

$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("D5").Value = "'230.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.46%  "

$ws.Range("D6").Value = "'0.5941"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.05%  "

$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("D8").Value = "'0.2768"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.75%  "

$ws.Range("D9").Value = "'0.06811"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.31%  "

$ws.Range("D10").Value = "'23.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.97%  "

$ws.Range("D11").Value = "'0.07520"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.01%  "

$ws.Range("D12").Value = "'1.802.37"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.44%  "

$ws.Range("D13").Value = "'4.692"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.85%  "

$ws.Range("D14").Value = "'0.6246"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.97%  "

$ws.Range("D15").Value = "'2.046.05"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.59%  "

$ws.Range("D16").Value = "'0.000009122"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -9.33%  "

$ws.Range("D17").Value = "'75.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.52%  "

$ws.Range("D18").Value = "'28.584.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.20%  "

$ws.Range("D19").Value = "'5.440"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -8.12%  "

$ws.Range("E20").Value = "  +0.13%  "

$ws.Range("D21").Value = "'210.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -8.34%  "

$ws.Range("E22").Value = "  -3.76%  "

$ws.Range("D23").Value = "'6.819"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.04%  "

$ws.Range("D24").Value = "'1.003"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.19%  "

$ws.Range("D25").Value = "'154.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.15%  "

$ws.Range("D26").Value = "'7.821"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.45%  "

$ws.Range("D27").Value = "'0.1274"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.54%  "

$ws.Range("D28").Value = "'16.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.86%  "

$ws.Range("D29").Value = "'1.442"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.20%  "

$ws.Range("D30").Value = "'0.06259"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.71%  "

$ws.Range("D31").Value = "'1.416"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.92%  "

$ws.Range("D32").Value = "'3.746"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.42%  "

$ws.Range("E33").Value = "  -3.07%  "

$ws.Range("D34").Value = "'1.704"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.30%  "

$ws.Range("E35").Value = "  -7.20%  "

$ws.Range("D36").Value = "'0.6342"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.05%  "

$ws.Range("D37").Value = "'2.504"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.24%  "

$ws.Range("D38").Value = "'2.713"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.90%  "

$ws.Range("D39").Value = "'0.01711"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.53%  "

$ws.Range("D40").Value = "'6.386"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.54%  "

$ws.Range("D41").Value = "'1.131.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.51%  "

$ws.Range("D42").Value = "'0.8628"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.26%  "

$ws.Range("E43").Value = "  +0.17%  "

$ws.Range("D44").Value = "'100.55"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.59%  "

$ws.Range("D45").Value = "'1.961.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.99%  "

$ws.Range("D46").Value = "'60.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.04%  "

$ws.Range("E47").Value = "  -5.56%  "

$ws.Range("D48").Value = "'1.572"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.36%  "

$ws.Range("D49").Value = "'8.348"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.82%  "

$ws.Range("D50").Value = "'0.4501"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.49%  "

$ws.Range("D51").Value = "'0.05444"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.77%  "
